# [ML] Change calculation of diff_mean to difference in means between the
# 'fake' and 'other' classes for each feature, where the data points for
# each feature are scaled to mean of 0 and standard deviation of 1.
#
# This updates the four recalculated statistic rows (diff mean value,
# diff mean percentage, p value, t value) and refreshes the banded row
# fill colors that are re-rolled whenever the report is regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated statistics for the affected features ---------------------

# Row 5: user_friends_count_per_day
$ws.Range("B5").Value = 1.843557631086175
$ws.Range("C5").Value = 0.5157993155584074
$ws.Range("D5").Value = 0.0008037934770308429
$ws.Range("E5").Value = 3.35322549038706

# Row 7: user_followers_count_per_day
$ws.Range("B7").Value = -450.2063900838606
$ws.Range("C7").Value = -0.392095559812901
$ws.Range("D7").Value = 0.0002528630914822416
$ws.Range("E7").Value = -3.661578251382192

# Row 16: created_at_weekday_sun_mon_tue
$ws.Range("B16").Value = -0.8667094347112938
$ws.Range("C16").Value = -0.2206850842391159
$ws.Range("D16").Value = 0.0006598264203504922
$ws.Range("E16").Value = -3.407571252368454

# Row 23: user_statuses_count_per_day
$ws.Range("B23").Value = 3.021404385478446
$ws.Range("C23").Value = 0.1939584297157631
$ws.Range("D23").Value = 0.002800949297166008
$ws.Range("E23").Value = 2.990038808446408

# --- Refreshed row banding colors --------------------------------------
# NOTE: multi-area ("A1:B1,C2:D2") ranges and Union() are not reliable in
# this COM engine (only the first area gets updated), so each contiguous
# block is set individually.

$ws.Range("A2:E2").Interior.Color = 0xE4A6BD
$ws.Range("A12:E12").Interior.Color = 0xE4A6BD
$ws.Range("A17:E18").Interior.Color = 0xE4A6BD
$ws.Range("A20:E21").Interior.Color = 0xE4A6BD

$ws.Range("A3:E4").Interior.Color = 0xCDE2FD
$ws.Range("A9:E10").Interior.Color = 0xCDE2FD
$ws.Range("A24:E25").Interior.Color = 0xCDE2FD
$ws.Range("A27:E29").Interior.Color = 0xCDE2FD
$ws.Range("A31:E31").Interior.Color = 0xCDE2FD
$ws.Range("A33:E33").Interior.Color = 0xCDE2FD
$ws.Range("A35:E35").Interior.Color = 0xCDE2FD

$ws.Range("A5:E5").Interior.Color = 0xE4CDA7

$ws.Range("A6:E6").Interior.Color = 0xF4FED5
$ws.Range("A32:E32").Interior.Color = 0xF4FED5
$ws.Range("A34:E34").Interior.Color = 0xF4FED5
$ws.Range("A36:E38").Interior.Color = 0xF4FED5

$ws.Range("A7:E8").Interior.Color = 0xB997B3

$ws.Range("A11:E11").Interior.Color = 0xCCA8EF
$ws.Range("A13:E13").Interior.Color = 0xCCA8EF

$ws.Range("A14:E15").Interior.Color = 0xA8DB99

$ws.Range("A16:E16").Interior.Color = 0xFFADDC

$ws.Range("A19:E19").Interior.Color = 0xE2DCDC

$ws.Range("A22:E22").Interior.Color = 0xCDDEAC

$ws.Range("A23:E23").Interior.Color = 0xF8D79D

$ws.Range("A26:E26").Interior.Color = 0x96DBB0
$ws.Range("A30:E30").Interior.Color = 0x96DBB0

$ws.Range("A39:E39").Interior.Color = 0xE1DBC3

Write-Host "edit applied"
